# Atualizacao de bases das ligas - swap Maccabi Netanya/Hapoel Beer Sheva
# and Hapoel Haifa/Hapoel TelAviv team-name references, plus corrected
# match records (row pairs 4/5, 10/11, 144/145, 222/223, 236/237 swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Israel Premier League")

$ws.Cells.Item(2,6).Value = 'Hapoel Haifa'  # F2
$ws.Cells.Item(4,2).Value = 6798412  # B4
$ws.Cells.Item(4,6).Value = 'Hapoel Hadera'  # F4
$ws.Cells.Item(4,7).Value = 3  # G4
$ws.Cells.Item(4,8).Value = 0  # H4
$ws.Cells.Item(4,9).Value = 2  # I4
$ws.Cells.Item(4,11).Value = 'H'  # K4
$ws.Cells.Item(4,12).Value = 1.3  # L4
$ws.Cells.Item(4,13).Value = 4.5  # M4
$ws.Cells.Item(4,14).Value = 8  # N4
$ws.Cells.Item(4,15).Value = 1.333  # O4
$ws.Cells.Item(4,16).Value = 4.333  # P4
$ws.Cells.Item(4,17).Value = 7  # Q4
$ws.Cells.Item(4,18).Value = -1.25  # R4
$ws.Cells.Item(4,19).Value = 1.825  # S4
$ws.Cells.Item(4,20).Value = 2.025  # T4
$ws.Cells.Item(4,21).Value = 2.75  # U4
$ws.Cells.Item(4,22).Value = 1.95  # V4
$ws.Cells.Item(4,23).Value = 1.9  # W4
$ws.Cells.Item(4,24).Value = 0.333  # X4
$ws.Cells.Item(4,25).Value = -1  # Y4
$ws.Cells.Item(4,27).Value = 0.825  # AA4
$ws.Cells.Item(4,28).Value = -1  # AB4
$ws.Cells.Item(4,29).Value = 0.475  # AC4
$ws.Cells.Item(4,30).Value = -0.5  # AD4
$ws.Cells.Item(5,2).Value = 6799822  # B5
$ws.Cells.Item(5,6).Value = 'Maccabi Bnei Raina'  # F5
$ws.Cells.Item(5,7).Value = 1  # G5
$ws.Cells.Item(5,8).Value = 1  # H5
$ws.Cells.Item(5,9).Value = 1  # I5
$ws.Cells.Item(5,11).Value = 'D'  # K5
$ws.Cells.Item(5,12).Value = 1.85  # L5
$ws.Cells.Item(5,13).Value = 3.3  # M5
$ws.Cells.Item(5,14).Value = 3.7  # N5
$ws.Cells.Item(5,15).Value = 1.65  # O5
$ws.Cells.Item(5,16).Value = 3.5  # P5
$ws.Cells.Item(5,17).Value = 4.5  # Q5
$ws.Cells.Item(5,18).Value = -0.75  # R5
$ws.Cells.Item(5,19).Value = 1.95  # S5
$ws.Cells.Item(5,20).Value = 1.9  # T5
$ws.Cells.Item(5,21).Value = 2.5  # U5
$ws.Cells.Item(5,22).Value = 2  # V5
$ws.Cells.Item(5,23).Value = 1.85  # W5
$ws.Cells.Item(5,24).Value = -1  # X5
$ws.Cells.Item(5,25).Value = 2.5  # Y5
$ws.Cells.Item(5,27).Value = -1  # AA5
$ws.Cells.Item(5,28).Value = 0.8999999999999999  # AB5
$ws.Cells.Item(5,29).Value = -1  # AC5
$ws.Cells.Item(5,30).Value = 0.8500000000000001  # AD5
$ws.Cells.Item(7,6).Value = 'Hapoel TelAviv'  # F7
$ws.Cells.Item(8,6).Value = 'Maccabi Netanya'  # F8
$ws.Cells.Item(10,2).Value = 6799829  # B10
$ws.Cells.Item(10,6).Value = 'Hapoel Beer Sheva'  # F10
$ws.Cells.Item(10,8).Value = 0  # H10
$ws.Cells.Item(10,9).Value = 1  # I10
$ws.Cells.Item(10,10).Value = 0  # J10
$ws.Cells.Item(10,11).Value = 'H'  # K10
$ws.Cells.Item(10,12).Value = 2.4  # L10
$ws.Cells.Item(10,13).Value = 3.4  # M10
$ws.Cells.Item(10,14).Value = 2.6  # N10
$ws.Cells.Item(10,15).Value = 2.625  # O10
$ws.Cells.Item(10,16).Value = 3.4  # P10
$ws.Cells.Item(10,17).Value = 2.375  # Q10
$ws.Cells.Item(10,18).Value = 0  # R10
$ws.Cells.Item(10,19).Value = 2  # S10
$ws.Cells.Item(10,20).Value = 1.85  # T10
$ws.Cells.Item(10,22).Value = 2  # V10
$ws.Cells.Item(10,23).Value = 1.85  # W10
$ws.Cells.Item(10,24).Value = 1.625  # X10
$ws.Cells.Item(10,25).Value = -1  # Y10
$ws.Cells.Item(10,27).Value = 1  # AA10
$ws.Cells.Item(10,28).Value = -1  # AB10
$ws.Cells.Item(10,29).Value = -1  # AC10
$ws.Cells.Item(10,30).Value = 0.8500000000000001  # AD10
$ws.Cells.Item(11,2).Value = 6799825  # B11
$ws.Cells.Item(11,6).Value = 'Maccabi Petach Tikva'  # F11
$ws.Cells.Item(11,8).Value = 2  # H11
$ws.Cells.Item(11,9).Value = 0  # I11
$ws.Cells.Item(11,10).Value = 1  # J11
$ws.Cells.Item(11,11).Value = 'D'  # K11
$ws.Cells.Item(11,12).Value = 1.8  # L11
$ws.Cells.Item(11,13).Value = 3.25  # M11
$ws.Cells.Item(11,14).Value = 4  # N11
$ws.Cells.Item(11,15).Value = 1.95  # O11
$ws.Cells.Item(11,16).Value = 3.2  # P11
$ws.Cells.Item(11,17).Value = 3.4  # Q11
$ws.Cells.Item(11,18).Value = -0.5  # R11
$ws.Cells.Item(11,19).Value = 2.1  # S11
$ws.Cells.Item(11,20).Value = 1.775  # T11
$ws.Cells.Item(11,22).Value = 2.025  # V11
$ws.Cells.Item(11,23).Value = 1.825  # W11
$ws.Cells.Item(11,24).Value = -1  # X11
$ws.Cells.Item(11,25).Value = 2.2  # Y11
$ws.Cells.Item(11,27).Value = -1  # AA11
$ws.Cells.Item(11,28).Value = 0.7749999999999999  # AB11
$ws.Cells.Item(11,29).Value = 1.025  # AC11
$ws.Cells.Item(11,30).Value = -1  # AD11
$ws.Cells.Item(16,5).Value = 'Hapoel TelAviv'  # E16
$ws.Cells.Item(19,5).Value = 'Hapoel Beer Sheva'  # E19
$ws.Cells.Item(21,5).Value = 'Maccabi Netanya'  # E21
$ws.Cells.Item(21,6).Value = 'Hapoel Haifa'  # F21
$ws.Cells.Item(22,6).Value = 'Maccabi Netanya'  # F22
$ws.Cells.Item(23,6).Value = 'Hapoel TelAviv'  # F23
$ws.Cells.Item(25,6).Value = 'Hapoel Beer Sheva'  # F25
$ws.Cells.Item(27,5).Value = 'Hapoel Haifa'  # E27
$ws.Cells.Item(32,5).Value = 'Hapoel TelAviv'  # E32
$ws.Cells.Item(33,5).Value = 'Maccabi Netanya'  # E33
$ws.Cells.Item(34,5).Value = 'Hapoel Beer Sheva'  # E34
$ws.Cells.Item(35,6).Value = 'Hapoel Haifa'  # F35
$ws.Cells.Item(37,5).Value = 'Hapoel Haifa'  # E37
$ws.Cells.Item(37,6).Value = 'Hapoel TelAviv'  # F37
$ws.Cells.Item(39,6).Value = 'Hapoel Beer Sheva'  # F39
$ws.Cells.Item(43,6).Value = 'Maccabi Netanya'  # F43
$ws.Cells.Item(44,6).Value = 'Hapoel Beer Sheva'  # F44
$ws.Cells.Item(46,5).Value = 'Hapoel TelAviv'  # E46
$ws.Cells.Item(47,6).Value = 'Hapoel Haifa'  # F47
$ws.Cells.Item(48,5).Value = 'Maccabi Netanya'  # E48
$ws.Cells.Item(52,5).Value = 'Hapoel Haifa'  # E52
$ws.Cells.Item(53,6).Value = 'Hapoel TelAviv'  # F53
$ws.Cells.Item(56,5).Value = 'Hapoel Beer Sheva'  # E56
$ws.Cells.Item(56,6).Value = 'Maccabi Netanya'  # F56
$ws.Cells.Item(57,6).Value = 'Hapoel Haifa'  # F57
$ws.Cells.Item(60,6).Value = 'Maccabi Netanya'  # F60
$ws.Cells.Item(62,6).Value = 'Hapoel Beer Sheva'  # F62
$ws.Cells.Item(63,5).Value = 'Hapoel TelAviv'  # E63
$ws.Cells.Item(66,5).Value = 'Hapoel Haifa'  # E66
$ws.Cells.Item(67,6).Value = 'Hapoel TelAviv'  # F67
$ws.Cells.Item(68,5).Value = 'Hapoel Beer Sheva'  # E68
$ws.Cells.Item(70,5).Value = 'Maccabi Netanya'  # E70
$ws.Cells.Item(71,5).Value = 'Hapoel Haifa'  # E71
$ws.Cells.Item(74,5).Value = 'Hapoel TelAviv'  # E74
$ws.Cells.Item(74,6).Value = 'Hapoel Beer Sheva'  # F74
$ws.Cells.Item(75,6).Value = 'Maccabi Netanya'  # F75
$ws.Cells.Item(78,5).Value = 'Hapoel Beer Sheva'  # E78
$ws.Cells.Item(79,6).Value = 'Hapoel Haifa'  # F79
$ws.Cells.Item(80,5).Value = 'Maccabi Netanya'  # E80
$ws.Cells.Item(80,6).Value = 'Hapoel TelAviv'  # F80
$ws.Cells.Item(86,6).Value = 'Hapoel Beer Sheva'  # F86
$ws.Cells.Item(88,6).Value = 'Maccabi Netanya'  # F88
$ws.Cells.Item(90,5).Value = 'Hapoel TelAviv'  # E90
$ws.Cells.Item(91,5).Value = 'Hapoel Haifa'  # E91
$ws.Cells.Item(92,6).Value = 'Hapoel Beer Sheva'  # F92
$ws.Cells.Item(93,6).Value = 'Maccabi Netanya'  # F93
$ws.Cells.Item(95,5).Value = 'Hapoel Haifa'  # E95
$ws.Cells.Item(98,5).Value = 'Hapoel TelAviv'  # E98
$ws.Cells.Item(99,5).Value = 'Maccabi Netanya'  # E99
$ws.Cells.Item(101,5).Value = 'Hapoel Beer Sheva'  # E101
$ws.Cells.Item(101,6).Value = 'Hapoel Haifa'  # F101
$ws.Cells.Item(104,6).Value = 'Hapoel TelAviv'  # F104
$ws.Cells.Item(109,6).Value = 'Hapoel TelAviv'  # F109
$ws.Cells.Item(110,6).Value = 'Hapoel Beer Sheva'  # F110
$ws.Cells.Item(111,5).Value = 'Hapoel Haifa'  # E111
$ws.Cells.Item(111,6).Value = 'Maccabi Netanya'  # F111
$ws.Cells.Item(114,5).Value = 'Hapoel TelAviv'  # E114
$ws.Cells.Item(115,5).Value = 'Hapoel Beer Sheva'  # E115
$ws.Cells.Item(116,5).Value = 'Maccabi Netanya'  # E116
$ws.Cells.Item(120,6).Value = 'Hapoel Haifa'  # F120
$ws.Cells.Item(122,6).Value = 'Hapoel TelAviv'  # F122
$ws.Cells.Item(123,6).Value = 'Maccabi Netanya'  # F123
$ws.Cells.Item(126,6).Value = 'Hapoel Beer Sheva'  # F126
$ws.Cells.Item(127,5).Value = 'Hapoel Haifa'  # E127
$ws.Cells.Item(129,5).Value = 'Hapoel Beer Sheva'  # E129
$ws.Cells.Item(132,5).Value = 'Hapoel TelAviv'  # E132
$ws.Cells.Item(132,6).Value = 'Hapoel Haifa'  # F132
$ws.Cells.Item(134,5).Value = 'Maccabi Netanya'  # E134
$ws.Cells.Item(135,5).Value = 'Hapoel Beer Sheva'  # E135
$ws.Cells.Item(136,6).Value = 'Maccabi Netanya'  # F136
$ws.Cells.Item(137,6).Value = 'Hapoel TelAviv'  # F137
$ws.Cells.Item(138,5).Value = 'Hapoel Haifa'  # E138
$ws.Cells.Item(143,6).Value = 'Hapoel Haifa'  # F143
$ws.Cells.Item(144,2).Value = 6799962  # B144
$ws.Cells.Item(144,5).Value = 'MS Ashdod'  # E144
$ws.Cells.Item(144,6).Value = 'Hapoel Petah Tikva'  # F144
$ws.Cells.Item(144,7).Value = 2  # G144
$ws.Cells.Item(144,9).Value = 2  # I144
$ws.Cells.Item(144,12).Value = 2.2  # L144
$ws.Cells.Item(144,13).Value = 3.1  # M144
$ws.Cells.Item(144,14).Value = 3.2  # N144
$ws.Cells.Item(144,15).Value = 2.2  # O144
$ws.Cells.Item(144,16).Value = 3.1  # P144
$ws.Cells.Item(144,17).Value = 3.2  # Q144
$ws.Cells.Item(144,18).Value = -0.25  # R144
$ws.Cells.Item(144,19).Value = 2  # S144
$ws.Cells.Item(144,20).Value = 1.85  # T144
$ws.Cells.Item(144,22).Value = 2  # V144
$ws.Cells.Item(144,23).Value = 1.85  # W144
$ws.Cells.Item(144,24).Value = 1.2  # X144
$ws.Cells.Item(144,27).Value = 1  # AA144
$ws.Cells.Item(144,29).Value = -0.5  # AC144
$ws.Cells.Item(144,30).Value = 0.425  # AD144
$ws.Cells.Item(145,2).Value = 6799960  # B145
$ws.Cells.Item(145,5).Value = 'Maccabi Petach Tikva'  # E145
$ws.Cells.Item(145,6).Value = 'Maccabi Bnei Raina'  # F145
$ws.Cells.Item(145,7).Value = 1  # G145
$ws.Cells.Item(145,9).Value = 0  # I145
$ws.Cells.Item(145,12).Value = 2.625  # L145
$ws.Cells.Item(145,13).Value = 3.25  # M145
$ws.Cells.Item(145,14).Value = 2.5  # N145
$ws.Cells.Item(145,15).Value = 2.8  # O145
$ws.Cells.Item(145,16).Value = 3.25  # P145
$ws.Cells.Item(145,17).Value = 2.375  # Q145
$ws.Cells.Item(145,18).Value = 0.25  # R145
$ws.Cells.Item(145,19).Value = 1.775  # S145
$ws.Cells.Item(145,20).Value = 2.1  # T145
$ws.Cells.Item(145,22).Value = 1.875  # V145
$ws.Cells.Item(145,23).Value = 1.975  # W145
$ws.Cells.Item(145,24).Value = 1.8  # X145
$ws.Cells.Item(145,27).Value = 0.7749999999999999  # AA145
$ws.Cells.Item(145,29).Value = -1  # AC145
$ws.Cells.Item(145,30).Value = 0.9750000000000001  # AD145
$ws.Cells.Item(147,5).Value = 'Maccabi Netanya'  # E147
$ws.Cells.Item(147,6).Value = 'Hapoel Beer Sheva'  # F147
$ws.Cells.Item(148,5).Value = 'Hapoel TelAviv'  # E148
$ws.Cells.Item(149,5).Value = 'Hapoel Haifa'  # E149
$ws.Cells.Item(150,5).Value = 'Maccabi Netanya'  # E150
$ws.Cells.Item(154,6).Value = 'Hapoel TelAviv'  # F154
$ws.Cells.Item(155,5).Value = 'Hapoel Beer Sheva'  # E155
$ws.Cells.Item(156,6).Value = 'Hapoel Haifa'  # F156
$ws.Cells.Item(157,5).Value = 'Hapoel TelAviv'  # E157
$ws.Cells.Item(159,6).Value = 'Hapoel Beer Sheva'  # F159
$ws.Cells.Item(160,6).Value = 'Maccabi Netanya'  # F160
$ws.Cells.Item(163,5).Value = 'Hapoel Beer Sheva'  # E163
$ws.Cells.Item(163,6).Value = 'Hapoel TelAviv'  # F163
$ws.Cells.Item(168,6).Value = 'Hapoel Haifa'  # F168
$ws.Cells.Item(169,5).Value = 'Maccabi Netanya'  # E169
$ws.Cells.Item(172,6).Value = 'Hapoel Beer Sheva'  # F172
$ws.Cells.Item(173,5).Value = 'Hapoel Haifa'  # E173
$ws.Cells.Item(176,5).Value = 'Hapoel TelAviv'  # E176
$ws.Cells.Item(176,6).Value = 'Maccabi Netanya'  # F176
$ws.Cells.Item(178,5).Value = 'Maccabi Netanya'  # E178
$ws.Cells.Item(179,5).Value = 'Hapoel Beer Sheva'  # E179
$ws.Cells.Item(182,6).Value = 'Hapoel TelAviv'  # F182
$ws.Cells.Item(183,6).Value = 'Hapoel Haifa'  # F183
$ws.Cells.Item(184,5).Value = 'Maccabi Netanya'  # E184
$ws.Cells.Item(184,6).Value = 'Hapoel TelAviv'  # F184
$ws.Cells.Item(188,5).Value = 'Hapoel Beer Sheva'  # E188
$ws.Cells.Item(190,6).Value = 'Hapoel Haifa'  # F190
$ws.Cells.Item(191,6).Value = 'Hapoel TelAviv'  # F191
$ws.Cells.Item(192,6).Value = 'Maccabi Netanya'  # F192
$ws.Cells.Item(196,5).Value = 'Hapoel Haifa'  # E196
$ws.Cells.Item(198,5).Value = 'Hapoel TelAviv'  # E198
$ws.Cells.Item(199,6).Value = 'Hapoel Beer Sheva'  # F199
$ws.Cells.Item(200,5).Value = 'Maccabi Netanya'  # E200
$ws.Cells.Item(201,6).Value = 'Hapoel TelAviv'  # F201
$ws.Cells.Item(202,5).Value = 'Hapoel Beer Sheva'  # E202
$ws.Cells.Item(204,6).Value = 'Maccabi Netanya'  # F204
$ws.Cells.Item(205,6).Value = 'Hapoel Haifa'  # F205
$ws.Cells.Item(210,6).Value = 'Hapoel Beer Sheva'  # F210
$ws.Cells.Item(212,5).Value = 'Hapoel Haifa'  # E212
$ws.Cells.Item(213,5).Value = 'Maccabi Netanya'  # E213
$ws.Cells.Item(214,5).Value = 'Hapoel TelAviv'  # E214
$ws.Cells.Item(215,5).Value = 'Hapoel Beer Sheva'  # E215
$ws.Cells.Item(215,6).Value = 'Hapoel Haifa'  # F215
$ws.Cells.Item(216,5).Value = 'Hapoel TelAviv'  # E216
$ws.Cells.Item(216,6).Value = 'Maccabi Netanya'  # F216
$ws.Cells.Item(222,2).Value = 8015674  # B222
$ws.Cells.Item(222,5).Value = 'Hapoel TelAviv'  # E222
$ws.Cells.Item(222,6).Value = 'Hapoel Bnei Sakhnin'  # F222
$ws.Cells.Item(222,7).Value = 1  # G222
$ws.Cells.Item(222,8).Value = 2  # H222
$ws.Cells.Item(222,9).Value = 0  # I222
$ws.Cells.Item(222,11).Value = 'A'  # K222
$ws.Cells.Item(222,12).Value = 1.95  # L222
$ws.Cells.Item(222,13).Value = 3.4  # M222
$ws.Cells.Item(222,14).Value = 3.75  # N222
$ws.Cells.Item(222,15).Value = 1.833  # O222
$ws.Cells.Item(222,16).Value = 3.5  # P222
$ws.Cells.Item(222,17).Value = 4.2  # Q222
$ws.Cells.Item(222,18).Value = -0.5  # R222
$ws.Cells.Item(222,19).Value = 1.825  # S222
$ws.Cells.Item(222,20).Value = 2.025  # T222
$ws.Cells.Item(222,22).Value = 2.025  # V222
$ws.Cells.Item(222,23).Value = 1.825  # W222
$ws.Cells.Item(222,24).Value = -1  # X222
$ws.Cells.Item(222,26).Value = 3.2  # Z222
$ws.Cells.Item(222,27).Value = -1  # AA222
$ws.Cells.Item(222,28).Value = 1.025  # AB222
$ws.Cells.Item(222,29).Value = 1.025  # AC222
$ws.Cells.Item(223,2).Value = 8015675  # B223
$ws.Cells.Item(223,5).Value = 'Maccabi Netanya'  # E223
$ws.Cells.Item(223,6).Value = 'Maccabi Bnei Raina'  # F223
$ws.Cells.Item(223,7).Value = 2  # G223
$ws.Cells.Item(223,8).Value = 1  # H223
$ws.Cells.Item(223,9).Value = 2  # I223
$ws.Cells.Item(223,11).Value = 'H'  # K223
$ws.Cells.Item(223,12).Value = 1.65  # L223
$ws.Cells.Item(223,13).Value = 3.6  # M223
$ws.Cells.Item(223,14).Value = 5.25  # N223
$ws.Cells.Item(223,15).Value = 1.5  # O223
$ws.Cells.Item(223,16).Value = 4  # P223
$ws.Cells.Item(223,17).Value = 6  # Q223
$ws.Cells.Item(223,18).Value = -1  # R223
$ws.Cells.Item(223,19).Value = 1.925  # S223
$ws.Cells.Item(223,20).Value = 1.925  # T223
$ws.Cells.Item(223,22).Value = 1.925  # V223
$ws.Cells.Item(223,23).Value = 1.925  # W223
$ws.Cells.Item(223,24).Value = 0.5  # X223
$ws.Cells.Item(223,26).Value = -1  # Z223
$ws.Cells.Item(223,27).Value = 0  # AA223
$ws.Cells.Item(223,28).Value = 0  # AB223
$ws.Cells.Item(223,29).Value = 0.925  # AC223
$ws.Cells.Item(225,5).Value = 'Hapoel Beer Sheva'  # E225
$ws.Cells.Item(227,5).Value = 'Hapoel Haifa'  # E227
$ws.Cells.Item(230,6).Value = 'Hapoel Haifa'  # F230
$ws.Cells.Item(231,6).Value = 'Maccabi Netanya'  # F231
$ws.Cells.Item(233,6).Value = 'Hapoel TelAviv'  # F233
$ws.Cells.Item(234,6).Value = 'Hapoel Beer Sheva'  # F234
$ws.Cells.Item(236,2).Value = 8016163  # B236
$ws.Cells.Item(236,5).Value = 'Hapoel TelAviv'  # E236
$ws.Cells.Item(236,6).Value = 'Maccabi Tel Aviv'  # F236
$ws.Cells.Item(236,7).Value = 0  # G236
$ws.Cells.Item(236,8).Value = 3  # H236
$ws.Cells.Item(236,9).Value = 0  # I236
$ws.Cells.Item(236,10).Value = 0  # J236
$ws.Cells.Item(236,12).Value = 3.4  # L236
$ws.Cells.Item(236,13).Value = 4  # M236
$ws.Cells.Item(236,14).Value = 1.75  # N236
$ws.Cells.Item(236,16).Value = 4.333  # P236
$ws.Cells.Item(236,17).Value = 1.571  # Q236
$ws.Cells.Item(236,18).Value = 1  # R236
$ws.Cells.Item(236,19).Value = 1.8  # S236
$ws.Cells.Item(236,20).Value = 2.05  # T236
$ws.Cells.Item(236,21).Value = 2.75  # U236
$ws.Cells.Item(236,26).Value = 0.571  # Z236
$ws.Cells.Item(236,28).Value = 1.05  # AB236
$ws.Cells.Item(236,29).Value = 0.4125  # AC236
$ws.Cells.Item(236,30).Value = -0.5  # AD236
$ws.Cells.Item(237,2).Value = 8016030  # B237
$ws.Cells.Item(237,5).Value = 'Maccabi Netanya'  # E237
$ws.Cells.Item(237,6).Value = 'Maccabi Haifa'  # F237
$ws.Cells.Item(237,7).Value = 1  # G237
$ws.Cells.Item(237,8).Value = 4  # H237
$ws.Cells.Item(237,9).Value = 1  # I237
$ws.Cells.Item(237,10).Value = 3  # J237
$ws.Cells.Item(237,12).Value = 2.625  # L237
$ws.Cells.Item(237,13).Value = 3.4  # M237
$ws.Cells.Item(237,14).Value = 2.25  # N237
$ws.Cells.Item(237,16).Value = 3.8  # P237
$ws.Cells.Item(237,17).Value = 1.727  # Q237
$ws.Cells.Item(237,18).Value = 0.75  # R237
$ws.Cells.Item(237,19).Value = 1.85  # S237
$ws.Cells.Item(237,20).Value = 2  # T237
$ws.Cells.Item(237,21).Value = 2.5  # U237
$ws.Cells.Item(237,26).Value = 0.7270000000000001  # Z237
$ws.Cells.Item(237,28).Value = 1  # AB237
$ws.Cells.Item(237,29).Value = 0.825  # AC237
$ws.Cells.Item(237,30).Value = -1  # AD237
$ws.Cells.Item(239,6).Value = 'Hapoel TelAviv'  # F239
$ws.Cells.Item(241,6).Value = 'Maccabi Netanya'  # F241
